$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 40/41: coin order swap (InjectiveProtocol <-> Aave) with updated price/volume ---
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'98.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.81%  "

# --- Price (column D) updates ---
$ws.Range("D2").Value = "35.597.90"
$ws.Range("D3").Value = "1.896.54"
$ws.Range("D5").Value = "'248.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Value = "'43.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.352"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.0741"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'13.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "2.171.82"
$ws.Range("D14").Value = "'0.730"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'4.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "1.895.98"
$ws.Range("D17").Value = "35.596.97"
$ws.Range("D18").Value = "'73.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "'247.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'12.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'165.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'8.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'18.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = "4.128.48"
$ws.Range("D31").Value = "'1.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'4.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.0582"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'4.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Value = "'0.851"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'2.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("D42").Value = "'0.0213"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'1.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "1.300.02"
$ws.Range("D45").Value = "'2.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0809"
$ws.Range("D46").Style = "Normal"
$ws.Range("D51").Value = "'43.48"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  +8.17%  "
$ws.Range("E9").Value = "  -5.29%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("E25").Value = "  -9.61%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +7.05%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -22.31%  "
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").Value = "  +7.53%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("E51").Value = "  -3.75%  "
